# Change the table style applied to the financial-documents table on
# slide 5 (the "B1- TYPES OF FINANCIAL DOCUMENTS" slide) from the
# deck's default/custom table style to the built-in
# "{489173CD-6D8B-442C-B066-8FC962375F51}" style.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)

foreach ($shape in $slide.Shapes) {
    if ($shape.HasTable) {
        $table = $shape.Table
        $table.ApplyStyle("{489173CD-6D8B-442C-B066-8FC962375F51}")
    }
}
